# Updating the pool names in ERM excel file
# Wrap each "E<n>-..." sheet name in parentheses around the number: "E(n)-..."

$wb = $excel.ActiveWorkbook

$renames = @(
    @{ old = "E1-Desktop Applications";      new = "E(1)-Desktop Applications" },
    @{ old = "E2-Electronic Messages";       new = "E(2)-Electronic Messages" },
    @{ old = "E3-Social Media";              new = "E(3)-Social Media" },
    @{ old = "E4-Cloud Services";            new = "E(4)-Cloud Services" },
    @{ old = "E5-Websites";                  new = "E(5)-Websites" },
    @{ old = "E6-Digital Media (Photo)";     new = "E(6)-Digital Media (Photo)" },
    @{ old = "E7-Digital Media (Audio)";     new = "E(7)-Digital Media (Audio)" },
    @{ old = "E8-Digital Media (Video)";     new = "E(8)-Digital Media (Video)" },
    @{ old = "E9-Databases";                 new = "E(9)-Databases" },
    @{ old = "E10-Shared Drives";            new = "E(10)-Shared Drives" },
    @{ old = "E11-Engineering Drawings";     new = "E(11)-Engineering Drawings" }
)

foreach ($r in $renames) {
    $wb.Worksheets($r.old).Name = $r.new
}

# The selection on "E(1)-Desktop Applications" moved from A49:Y49 to a single
# cell C49 (scrolled back up to the top of the sheet as well).
$ws1 = $wb.Worksheets("E(1)-Desktop Applications")
$ws1.Activate()
$ws1.Range("C49").Select()

# "E(11)-Engineering Drawings" becomes the active/selected tab when the
# workbook is saved (was "ERM Vendors Master File" before).
$wsLast = $wb.Worksheets("E(11)-Engineering Drawings")
$wsLast.Activate()
